$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Y of" value used in the sum (D6: 2 -> 3)
$ws.Range("D6").Value = 3

# Replace the literal result with a formula that actually sums D5 and D6
$ws.Range("D11").Formula = "=D5+D6"

# Rename the worksheet and update the matching title cell to match
$ws.Name = "Sum 1 and 3"
$ws.Range("B1").Value = "Sum 1 and 3"

# Move the selection to B1, matching the saved view state
$ws.Range("B1").Select()
